# Reformat the sole shared string (JSON-pretty-print it), and move it
# from A2 up to A1 - replacing the old numeric placeholder cell (which
# carried the bold/bordered/centered "header" style) so A1 ends up plain
# (un-styled) and holds the question-set text, with A2 removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = 'questions = [
    {
        "title": "You are a marketing analyst in an e-shop. You see that the sales decreased last week, which is unusual for this period, so you look for what on the website might be causing the decrease in sales. The only analytical tool your company uses is Google Analytics. What Google Analytics data will you analyze to find the weak spot?",
        "ques_type": 2,
        "options": [
            "Shopping behavior: consumers\u2019 path on the website and inside the cart",
            "Audience: demographics data and users\u2019 interests",
            "Acquisition: paid and owned channels\u2019 effectiveness",
            "Behavior flow: actions from the landing pages through exit pages"
        ],
        "score": "Shopping behavior: consumers\u2019 path on the website and inside the cart"
    },
    {
        "title": "You are just ending your first year at a company, where you are responsible for digital marketing analytics. The Head of Marketing requests that you create an annual report. However, the company has not created previous reports or collected data for previous years except for the basic analytical tools, such as Google Analytics or Hotjar.  What metrics can you use to compare Search Engine Optimization effectiveness year to year based on the data you have?",
        "ques_type": 2,
        "options": [
            "Position of the website in search engines for the main keywords",
            "Number of keywords the website appears on first 100 positions for",
            "Organic search traffic to the website",
            "Number of backlinks for the website"
        ],
        "score": "Organic search traffic to the website"
    },
    {
        "title": "As your company''s marketing analyst, you are asked to provide information on the age structure of your customers. Here is the data you see:18-25 years old: 380 people26-35 years old: 100 people36-45 years old: 120 people46-55 years old: 90 people56+ years old: 10 peopleUnknown: 300 people What is the age of your core audience?",
        "ques_type": 2,
        "options": [
            "The core audience is under 25 years old.",
            "The core audience is 18-25 years old.",
            "The core audience is 18-35 years old.",
            "The core audience is 18-55 years old."
        ],
        "score": "The core audience is under 25 years old."
    },
    {
        "title": "Making the report on your company website''s organic traffic, you see that your e-shop has high traffic, high bounce rate, low average duration of a session, and low average number of pages per session.  What would you recommend checking and optimizing to improve this situation?",
        "ques_type": 15,
        "options": [
            "Usability of the website",
            "Structure of the cart",
            "Effectiveness of advertising",
            "Content in social media",
            "Content of emails",
            "Content of the website"
        ],
        "score": [
            "Usability of the website",
            "Content in social media",
            "Content of the website"
        ]
    }
]'

# A2 (the old shared-string cell) is no longer needed - the sheet will
# only contain a single cell once this edit lands.
$ws.Range("A2").ClearContents()

# Drop the bold/border/center-aligned style that used to live on A1 (it
# held a placeholder 0) before writing the real text into it.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# Writing text with embedded line breaks makes the host auto-size the row;
# put the row back to the sheet's default automatic height.
$ws.Rows(1).EntireRow.AutoFit()
